# "Review Coding Phase & Update Code"
# Fills in the Coding Phase Defects review table (rows 10-16) with the
# reviewed defect entries found while reviewing the coding phase / code.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coding Phase Defects")

# --- Row 10: blank-field validation message (AddPartController & friends) ---
$ws.Range("C10").Value = 'C08'
$ws.Range("D10").Value = 'AddPartController.java, 166; AddProductController.java, 212; ModifyPartController.java, 193; ModifyProductController.java, 236'
$ws.Range("E10").Value = 'Se afișează următorul mesaj: "Form contains blank field." în momentul în care se încearcă salvarea unor date greșite (tipul diferă, în loc de int se introduce string). Ar trebui afișat un mesaj corespunzător problemei.'

# --- Row 11: no part selected on modify ---
$ws.Range("C11").Value = 'C06'
$ws.Range("D11").Value = 'ModifyPartController.java, 139'
$ws.Range("E11").Value = 'Daca nu se selecteză o piesă pentru modificare, în interfața cu utilizatorul ar trebui să apară un mesaj prin care să fie atenționat că este necesar să selecteze o piesă. În momentul de față, nu apare nimic pe interfață, dar se aruncă o eroare (java.lang.IndexOutOfBoundsException), deoarece nu se verifică dacă a fost selectată o piesă și se încearcă apelarea metodei getPartId() pe o piesă cu valoarea null.'

# --- Row 12: no part selected on delete ---
$ws.Range("C12").Value = 'C06'
$ws.Range("D12").Value = 'ModifyPartController.java, 139'
$ws.Range("E12").Value = 'Daca nu se selecteză o piesă pentru ștergere în interfața cu utilizatorul ar trebui să apară un mesaj prin care să fie atenționat că este necesar să selecteze o piesă. În momentul de față, nu apare nimic pe interfață, dar se aruncă o eroare (java.lang.NullPointerException), deoarece nu se verifică dacă a fost selectată o piesă și se încearcă apelarea metodei getName() pe o piesă cu valoarea null.'

# --- Row 13: inStock compared to 1 instead of 0 ---
$ws.Range("C13").Value = 'C01'
$ws.Range("D13").Value = 'Part.java, 93'
$ws.Range("E13").Value = 'Variabila inStock e comparată cu 1 în loc de 0'

# --- Row 14: missing resource file -> NullPointerException ---
$ws.Range("C14").Value = 'C01'
$ws.Range("D14").Value = 'InventoryRepository.java, 26, 73, 121'
$ws.Range("E14").Value = 'Dacă fisierul cu date nu există în resurse, aplicația va arunca NullPointerException. Ar trebui reținut rezultatul apelării getResource, și comparat cu null, dacă este cazul aruncată o eroare adecvată '

# --- Row 15: no object selected on modify/delete (MainScreenController) ---
$ws.Range("C15").Value = 'C06'
$ws.Range("D15").Value = 'MainScreenController.java, 210,223,158,128'
$ws.Range("E15").Value = 'La operațiile de modificare și ștergere, în interfața cu utilizatorul ar trebui să apară un mesaj prin care să fie atenționat că este necesar să selecteze un obiect, in momentul de față, dacă un obiect nu e selectat, valoarea va fi null si se continua cu null în funcțiile urmatoare'

# --- Row 16: String.format uses "{}" instead of "%s" ---
$ws.Range("C16").Value = 'C04'
$ws.Range("D16").Value = 'MainScreenController.java, 185,187'
$ws.Range("E16").Value = 'Pentru metoda String.format se folosește "{}" in loc de "%s" '

# Row heights to fit the new (wrapped) text
$ws.Rows.Item(10).RowHeight = 120
$ws.Rows.Item(11).RowHeight = 150
$ws.Rows.Item(12).RowHeight = 150
$ws.Rows.Item(13).RowHeight = 15
$ws.Rows.Item(14).RowHeight = 75
$ws.Rows.Item(15).RowHeight = 105
$ws.Rows.Item(16).RowHeight = 30

# Column widths so the new "File, Line" / comment text is readable
$ws.Columns.Item(4).ColumnWidth = 25
$ws.Columns.Item(6).ColumnWidth = 37.71

# Stray formatting artifact next to the row-11 comment (matches source edit)
$ws.Range("F11").WrapText = $true

# E15/E16 picked up from a paste without the table border
$ws.Range("E15").Borders.LineStyle = 0
$ws.Range("E16").Borders.LineStyle = 0
$ws.Range("E15").WrapText = $true
$ws.Range("E16").WrapText = $true

# Leave the view scrolled down to the newly-filled rows, like the author did
$ws.Activate()
$ws.Range("E53").Select()
